$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> { D = new price text; E = new volume text }
$data = @{
    2 = @{ D = "60.292.36"; E = "  +5.42%  " }
    3 = @{ D = "2.616.14"; E = "  +7.94%  " }
    4 = @{ D = "1.01"; E = "  +0.44%  " }
    5 = @{ D = "506.35"; E = "  +3.67%  " }
    6 = @{ D = "156.56"; E = "  +1.59%  " }
    7 = @{ D = "0.998"; E = "  +0.26%  " }
    8 = @{ D = "0.587"; E = "  -5.09%  " }
    9 = @{ D = "2.630.34"; E = "  +7.56%  " }
    10 = @{ E = "  +5.26%  " }
    11 = @{ D = "0.105"; E = "  +4.37%  " }
    12 = @{ E = "  +2.82%  " }
    13 = @{ E = "  +0.78%  " }
    14 = @{ D = "3.080.33"; E = "  +8.03%  " }
    15 = @{ D = "60.407.67"; E = "  +5.60%  " }
    16 = @{ D = "21.69"; E = "  +5.28%  " }
    17 = @{ D = "0.0000140"; E = "  +4.95%  " }
    18 = @{ D = "2.627.54"; E = "  +7.87%  " }
    19 = @{ D = "4.77"; E = "  +3.01%  " }
    20 = @{ D = "344.49"; E = "  +5.95%  " }
    21 = @{ D = "10.42"; E = "  +4.20%  " }
    22 = @{ E = "  +3.79%  " }
    23 = @{ D = "0.999"; E = "  +0.02%  " }
    24 = @{ D = "60.05"; E = "  +3.81%  " }
    25 = @{ D = "0.423"; E = "  +5.41%  " }
    26 = @{ E = "  +2.97%  " }
    27 = @{ D = "0.994"; E = "  -0.29%  " }
    28 = @{ D = "0.0₃0855"; E = "  +9.16%  " }
    29 = @{ D = "7.55"; E = "  +3.64%  " }
    30 = @{ D = "1.00"; E = "  +0.19%  " }
    31 = @{ D = "156.67" }
    32 = @{ D = "19.36"; E = "  +3.42%  " }
    33 = @{ E = "  +3.35%  " }
    34 = @{ D = "5.71"; E = "  +7.69%  " }
    35 = @{ E = "  +6.08%  " }
    36 = @{ E = "  +4.70%  " }
    37 = @{ D = "307.75"; E = "  +7.74%  " }
    38 = @{ D = "0.848"; E = "  +3.91%  " }
    39 = @{ E = "  +7.36%  " }
    40 = @{ E = "  +6.96%  " }
    41 = @{ D = "0.831"; E = "  +26.83%  " }
    42 = @{ D = "35.58"; E = "  +4.57%  " }
    43 = @{ D = "0.627"; E = "  +4.65%  " }
    44 = @{ D = "0.0571"; E = "  +7.55%  " }
    45 = @{ D = "0.101"; E = "  -1.27%  " }
    46 = @{ D = "0.993"; E = "  -0.12%  " }
    47 = @{ E = "  +12.65%  " }
    48 = @{ D = "4.86"; E = "  +6.67%  " }
    49 = @{ D = "0.0236"; E = "  +3.73%  " }
    50 = @{ D = "2.048.25"; E = "  +7.82%  " }
    51 = @{ D = "10.29"; E = "  +0.74%  " }
}

# Rows whose Price text looks like a plain number (e.g. "1.01", "0.0000140")
# need the cell pre-formatted as Text, otherwise Excel auto-converts the
# assignment into a numeric value and the literal formatting (trailing
# zeros, exponent form, etc.) is lost.
$textFormatRows = @(4, 5, 6, 7, 8, 11, 16, 17, 19, 20, 21, 23, 24, 25, 27, 29, 30, 31, 32, 34, 37, 38, 41, 42, 43, 44, 45, 46, 48, 49, 51)

foreach ($row in ($data.Keys | Sort-Object)) {
    $vals = $data[$row]
    if ($vals.ContainsKey("D")) {
        $dCell = $ws.Cells.Item($row, 4)
        if ($textFormatRows -contains $row) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $vals["E"]
    }
}
